$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference a known General/style-0, text-holding cell (D4 is untouched by this edit)
# to restore style after forcing text-format on numeric-looking Price values.
$refStyle = $ws.Range("D4").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.253.52'
$ws.Range("D2").Style = $refStyle
$ws.Range("E2").Value = '  +1.98%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.172.83'
$ws.Range("D3").Style = $refStyle
$ws.Range("E3").Value = '  +3.89%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.46'
$ws.Range("D5").Style = $refStyle
$ws.Range("E5").Value = '  +2.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.43'
$ws.Range("D6").Style = $refStyle
$ws.Range("E6").Value = '  +6.18%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.167.92'
$ws.Range("D8").Style = $refStyle
$ws.Range("E8").Value = '  +3.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("D9").Style = $refStyle
$ws.Range("E9").Value = '  +2.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.164'
$ws.Range("D10").Style = $refStyle
$ws.Range("E10").Value = '  +5.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.19'
$ws.Range("D11").Style = $refStyle
$ws.Range("E11").Value = '  +2.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.507'
$ws.Range("D12").Style = $refStyle
$ws.Range("E12").Value = '  +6.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000276'
$ws.Range("D13").Style = $refStyle
$ws.Range("E13").Value = '  +19.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.19'
$ws.Range("D14").Style = $refStyle
$ws.Range("E14").Value = '  +9.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.691.42'
$ws.Range("D15").Style = $refStyle
$ws.Range("E15").Value = '  +4.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.268.31'
$ws.Range("D16").Style = $refStyle
$ws.Range("E16").Value = '  +1.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.23'
$ws.Range("D17").Style = $refStyle
$ws.Range("E17").Value = '  +7.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.178.86'
$ws.Range("D18").Style = $refStyle
$ws.Range("E18").Value = '  +4.02%  '
$ws.Range("E19").Value = '  +1.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '512.40'
$ws.Range("D20").Style = $refStyle
$ws.Range("E20").Value = '  +7.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.95'
$ws.Range("D21").Style = $refStyle
$ws.Range("E21").Value = '  +6.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.06'
$ws.Range("D22").Style = $refStyle
$ws.Range("E22").Value = '  +13.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.738'
$ws.Range("D23").Style = $refStyle
$ws.Range("E23").Value = '  +9.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.87'
$ws.Range("D24").Style = $refStyle
$ws.Range("E24").Value = '  +4.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.97'
$ws.Range("D25").Style = $refStyle
$ws.Range("E25").Value = '  +3.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = $refStyle
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("E27").Value = '  +15.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.91'
$ws.Range("D28").Style = $refStyle
$ws.Range("E28").Value = '  +4.36%  '
$ws.Range("E29").Value = '  +9.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '28.13'
$ws.Range("D30").Style = $refStyle
$ws.Range("E30").Value = '  +7.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.78'
$ws.Range("D31").Style = $refStyle
$ws.Range("E31").Value = '  +14.56%  '
$ws.Range("E32").Value = '  +7.96%  '
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.33'
$ws.Range("D34").Style = $refStyle
$ws.Range("E34").Value = '  +13.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.68'
$ws.Range("D35").Style = $refStyle
$ws.Range("E35").Value = '  +8.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.65'
$ws.Range("D36").Style = $refStyle
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '477.81'
$ws.Range("D37").Style = $refStyle
$ws.Range("E37").Value = '  +8.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.16'
$ws.Range("D38").Style = $refStyle
$ws.Range("E38").Value = '  +12.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0879'
$ws.Range("D39").Style = $refStyle
$ws.Range("E39").Value = '  +9.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0422'
$ws.Range("D40").Style = $refStyle
$ws.Range("E40").Value = '  +4.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.143.71'
$ws.Range("D41").Style = $refStyle
$ws.Range("E41").Value = '  +6.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.64'
$ws.Range("D42").Style = $refStyle
$ws.Range("E42").Value = '  +5.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.122'
$ws.Range("D43").Style = $refStyle
$ws.Range("E43").Value = '  +7.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.52'
$ws.Range("D44").Style = $refStyle
$ws.Range("E44").Value = '  +17.21%  '
$ws.Range("E45").Value = '  +11.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.25'
$ws.Range("D46").Style = $refStyle
$ws.Range("E46").Value = '  +6.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0587'
$ws.Range("D47").Style = $refStyle
$ws.Range("E47").Value = '  +14.35%  '
$ws.Range("E49").Value = '  +2.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.33'
$ws.Range("D50").Style = $refStyle
$ws.Range("E50").Value = '  +12.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '123.74'
$ws.Range("D51").Style = $refStyle
$ws.Range("E51").Value = '  +5.70%  '
